$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mu = [char]0x03BC   # micro sign, U+03BC

# Replace the ascii "u" micro-prefix abbreviation with the proper "μ" sign
# in the capacitor value column (B) for the affected rows.
$ws.Range("B2").Value2  = "0.33$mu"    # C1
$ws.Range("B3").Value2  = "33$mu"      # C2
$ws.Range("B4").Value2  = "0.1$mu"     # C3
$ws.Range("B7").Value2  = "4.7$mu"     # C6
$ws.Range("B9").Value2  = "0.01$mu"    # C8
$ws.Range("B10").Value2 = "2.2$mu"     # C9
$ws.Range("B23").Value2 = "0.001$mu"   # C22
$ws.Range("B26").Value2 = "1$mu"       # C25
$ws.Range("B28").Value2 = "0.18$mu"    # C27
$ws.Range("B34").Value2 = "10$mu"      # C33
$ws.Range("B37").Value2 = "150$mu"     # C36

# P3 value: "JST XH header" -> "XH header"
$ws.Range("B48").Value2 = "XH header"

# Update the saved view state: scroll position and active cell/selection.
$ws.Application.ActiveWindow.ScrollRow = 23
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B51").Select()
